# Auto-generated edit script
# Sets column B (退服) values for specific rows on both sheets, per commit diff.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$sheet1Updates = @(
@{Row=7; Value=1}, @{Row=10; Value=1}, @{Row=41; Value=1}, @{Row=114; Value=1}, @{Row=116; Value=1}, @{Row=140; Value=1}, @{Row=224; Value=1}, @{Row=253; Value=1}, @{Row=321; Value=1}, @{Row=339; Value=1}, @{Row=343; Value=1}, @{Row=345; Value=1}, @{Row=348; Value=1}, @{Row=433; Value=0}, @{Row=456; Value=0}, @{Row=461; Value=0}, @{Row=462; Value=1}, @{Row=506; Value=1}, @{Row=540; Value=1}, @{Row=590; Value=1}, @{Row=596; Value=1}, @{Row=669; Value=1}, @{Row=699; Value=1}, @{Row=700; Value=1}, @{Row=733; Value=1}, @{Row=741; Value=1}, @{Row=754; Value=0}, @{Row=758; Value=1}, @{Row=783; Value=0}, @{Row=824; Value=1}, @{Row=841; Value=1}, @{Row=845; Value=1}, @{Row=865; Value=1}, @{Row=874; Value=0}, @{Row=916; Value=1}, @{Row=948; Value=1}, @{Row=949; Value=1}, @{Row=966; Value=1}, @{Row=974; Value=0}, @{Row=978; Value=1}, @{Row=997; Value=0}
)

$sheet2Updates = @(
@{Row=50; Value=1}, @{Row=83; Value=1}, @{Row=99; Value=1}, @{Row=176; Value=0}, @{Row=308; Value=0}, @{Row=334; Value=0}, @{Row=393; Value=0}, @{Row=397; Value=1}, @{Row=464; Value=0}, @{Row=471; Value=0}, @{Row=484; Value=0}, @{Row=522; Value=0}, @{Row=531; Value=0}, @{Row=646; Value=0}, @{Row=647; Value=0}, @{Row=683; Value=0}, @{Row=741; Value=1}, @{Row=770; Value=0}, @{Row=799; Value=0}, @{Row=827; Value=0}, @{Row=973; Value=1}
)

foreach ($u in $sheet1Updates) {
    $ws1.Cells.Item($u.Row, 2).Value = $u.Value
}

foreach ($u in $sheet2Updates) {
    $ws2.Cells.Item($u.Row, 2).Value = $u.Value
}

Write-Host "Applied $($sheet1Updates.Count) updates to sheet '$($ws1.Name)' and $($sheet2Updates.Count) updates to sheet '$($ws2.Name)'."
